$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.405.08"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.008.26"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'256.92"
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'55.59"
$ws.Range("E8").Value = "  -4.19%  "
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").Value = "'0.0764"
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "2.308.99"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "'14.21"
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("D14").Value = "'21.79"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "'0.773"
$ws.Range("E15").Value = "  -6.80%  "
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").Value = "2.031.73"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "37.273.61"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "'69.72"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "0.0₃0827"
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").Value = "'234.24"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'2.54"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").Value = "'164.74"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'8.85"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.131"
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.40"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("D33").Value = "'0.0624"
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("D34").Value = "'4.36"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("D35").Value = "'2.39"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D36").Value = "'3.42"
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'5.35"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "1.443.01"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").Value = "'0.0917"
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("D44").Value = "'0.0208"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'15.78"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'89.29"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "'6.80"
$ws.Range("E49").Value = "  -9.08%  "
$ws.Range("D50").Value = "2.198.85"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'1.90"
$ws.Range("E51").Value = "  -8.12%  "
